# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the newly generated data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row number => new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 1079
    3  = 372
    4  = 1480
    5  = 8710
    6  = 87
    10 = 150
    11 = 12
    12 = 3559
    14 = 362
    16 = 1155
    18 = 1114
    20 = 199
    21 = 2306
    22 = 53
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (All types) - row number => new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 1079
    3  = 372
    4  = 1480
    5  = 8710
    6  = 87
    10 = 150
    11 = 12
    12 = 3559
    14 = 362
    16 = 1155
    18 = 1114
    20 = 199
    21 = 2306
    23 = 53
}
foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allTypesUpdates[$row]
}
